$d = $word.ActiveDocument

$bodyXml = @'
<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:cs/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc379225272"/><w:r><w:t>URS-SRS</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>URS-01: The user has to register themselves as customer</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-01: The system provide an interface for registration which receive name address and their picture for the registration. </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-02: The system can check the name format. The name format should be only a character.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-03: The system can check the address. The address format should be a character and number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-04: The system can check the picture that user upload. The picture should be an image format (for example: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jpg</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-05: The system can check the username format. The username format should be character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-06: The system can check the password format. The password format should be character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-07: The system shall retrieve the existing data from the customer username</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-08: The system shall add a new customer information to database using username, password, name, address and their picture.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-09: </w:t></w:r><w:r><w:t>The system shall provide the successful registration page</w:t></w:r><w:r><w:t xml:space="preserve"> after registration is complete</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-10: </w:t></w:r><w:r><w:t>The system shall display the error message “The name must be only characters”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-11: </w:t></w:r><w:r><w:t>The system shall display the error message “</w:t></w:r><w:r><w:t>The picture that must be image format</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-12: </w:t></w:r><w:r><w:t>The system shall dis</w:t></w:r><w:r><w:t xml:space="preserve">play the error message “The username must be </w:t></w:r><w:r><w:t>characters</w:t></w:r><w:r><w:t xml:space="preserve"> with number</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-13: </w:t></w:r><w:r><w:t>The system shall dis</w:t></w:r><w:r><w:t xml:space="preserve">play the error message “The password </w:t></w:r><w:r><w:t>must</w:t></w:r><w:r><w:t xml:space="preserve"> be </w:t></w:r><w:r><w:t>characters</w:t></w:r><w:r><w:t xml:space="preserve"> with number</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="1" w:name="_Toc379225273"/><w:r><w:t>URS-02: The user log in as customer</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-14: </w:t></w:r><w:r><w:t>The system provides the login user interface which receive the user name and password</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-05</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t xml:space="preserve">The system can check the username format. The username format should be </w:t></w:r><w:r><w:t>character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-06</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> The system can check the pa</w:t></w:r><w:r><w:t>ssword format.</w:t></w:r><w:r><w:t xml:space="preserve"> The password must contain the </w:t></w:r><w:r><w:t>character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>15</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system validates the username and password</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>SRS-</w:t></w:r><w:r><w:t>16</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall pro</w:t></w:r><w:r><w:t>vide the successful login</w:t></w:r><w:r><w:t xml:space="preserve"> page.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-17: The system shall provide catalog page that consist with product name, product description and product picture.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-18</w:t></w:r><w:r><w:t>: The system shall show the customer status with a customer picture on the top-right of the user interface to confirm his identity.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="2" w:name="_Toc379225274"/><w:r><w:t>URS-03: The customer update their information</w:t></w:r><w:bookmarkEnd w:id="2"/></w:p><w:p><w:r><w:tab/><w:t>SRS-19</w:t></w:r><w:r><w:t>: The system provide an interface for registration which receive name address and their picture for the registration.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-02</w:t></w:r><w:r><w:t>: The system can check the name format. The name format should be only a character.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-03</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system can check the address. The address format should be a character and number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-04</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t xml:space="preserve"> The system can check the picture that user update. The picture should be </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> image format (for example: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jpg</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-20</w:t></w:r><w:r><w:t>: The system shall update customer information to database using name, address and their picture.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-21</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall provide the successful registration page</w:t></w:r><w:r><w:t xml:space="preserve"> after customer information updating was complete</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-10</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall display the error message “The name must be only characters”.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-11</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall display the error message “</w:t></w:r><w:r><w:t>The picture that must be image format</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="3" w:name="_Toc379225275"/><w:r><w:t>URS-04: The customer search product.</w:t></w:r><w:bookmarkEnd w:id="3"/></w:p><w:p><w:r><w:tab/><w:t>SRS-22</w:t></w:r><w:r><w:t>: The system shall provide search product page which receive product name for searching.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-23</w:t></w:r><w:r><w:t xml:space="preserve">: The system shall retrieve the existing data from database follow product name that user input. </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-24</w:t></w:r><w:r><w:t>: The system shall provide detail</w:t></w:r><w:r><w:t xml:space="preserve"> page that consist with product name, product description</w:t></w:r><w:r><w:t>, item quantity</w:t></w:r><w:r><w:t xml:space="preserve"> and product picture.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>URS -05: The customer add product to a cart</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">SRS-17: </w:t></w:r><w:r><w:t>The system shall provide catalog page that consist with product name, product description and product picture.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-24: The system shall provide detail page that consist with product name, product description, item quantity and product picture.</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-25: The system shall provide the interface for customer can select the product to a cart.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-26: The system shall provide the interface for customer can select the product quantity to a cart.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-27: The system save cart history into a database using customer</w:t></w:r><w:r><w:t xml:space="preserve"> ID, product ID and status</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">SRS-28: The system update status at cart history </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:bookmarkStart w:id="4" w:name="_Toc379225277"/><w:r><w:t>SRS-29</w:t></w:r><w:r><w:t>: Th</w:t></w:r><w:r><w:t>e system shall display the alert</w:t></w:r><w:r><w:t xml:space="preserve"> message “</w:t></w:r><w:r><w:t>The product has been add to a cart</w:t></w:r><w:r><w:t>”.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>URS-06: The customer check out for purchase the product.</w:t></w:r><w:bookmarkEnd w:id="4"/></w:p><w:p><w:r><w:tab/><w:t>SRS-30</w:t></w:r><w:r><w:t>: The system shall provide interface for check out for purchase the product.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-31</w:t></w:r><w:r><w:t>: The system shall provide the product list that consist product name, product picture, price, quantity and total price.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-32</w:t></w:r><w:r><w:t>: The system shall provide payment method interface.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-33</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall pro</w:t></w:r><w:r><w:t>vide the successful</w:t></w:r><w:r><w:t xml:space="preserve"> purchasing</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>page</w:t></w:r><w:r><w:t xml:space="preserve"> after purchase is complete</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-34: The system shall update status at cart history in a database.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-35</w:t></w:r><w:r><w:t>: The system shall add a new shopping history to database.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="5" w:name="_Toc379225278"/><w:r><w:t>URS-07: The customer view their history.</w:t></w:r><w:bookmarkEnd w:id="5"/></w:p><w:p><w:r><w:tab/><w:t>SRS-</w:t></w:r><w:r><w:t>36</w:t></w:r><w:r><w:t>: The system shall provide interface for view a customer history.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-37</w:t></w:r><w:r><w:t>: The system query history list page that consist with date and history id.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-38</w:t></w:r><w:r><w:t>: The system query the history detail that consist the list of product that customer bought and total price.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>URS-08: The user log in as administrator</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-14</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system provides the login user interface which receive the user name and password</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-05: </w:t></w:r><w:r><w:t xml:space="preserve">The system can check the username format. The username format should be </w:t></w:r><w:r><w:t>character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-06:</w:t></w:r><w:r><w:t xml:space="preserve"> The system can check the pa</w:t></w:r><w:r><w:t>ssword format.</w:t></w:r><w:r><w:t xml:space="preserve"> The password must contain the </w:t></w:r><w:r><w:t>character with number.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>15</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system validates the username and password</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t xml:space="preserve">SRS-16: </w:t></w:r><w:r><w:t>The system shall pro</w:t></w:r><w:r><w:t>vide the successful login</w:t></w:r><w:r><w:t xml:space="preserve"> page.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>SRS-39</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The system shall provide administrator page after login successful.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="6" w:name="_Toc379225280"/><w:r><w:t>URS-09: The administrator add a new product.</w:t></w:r><w:bookmarkEnd w:id="6"/></w:p><w:p><w:r><w:tab/><w:t>SRS-40</w:t></w:r><w:r><w:t>: The system shall provide the interface for add new product which receive product name, description and product picture.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-41</w:t></w:r><w:r><w:t>: The system can check the product title format. The title should be character.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">SRS-04: The system can check the picture that user upload. The picture should be an image format (for example: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jpg</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-42</w:t></w:r><w:r><w:t>: The system shall retrieve the existing data from the product name.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-43</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>The</w:t></w:r><w:r><w:t xml:space="preserve"> system shall add a new product</w:t></w:r><w:r><w:t xml:space="preserve"> inform</w:t></w:r><w:r><w:t>ation to database using product name, description and product picture</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-44</w:t></w:r><w:r><w:t>: The</w:t></w:r><w:r><w:t xml:space="preserve"> system shall pro</w:t></w:r><w:r><w:t>vide the successful add new product</w:t></w:r><w:r><w:t xml:space="preserve"> page.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-45</w:t></w:r><w:r><w:t>: The system shall display the error message “This product is existing in database”.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:bookmarkStart w:id="7" w:name="_Toc379225281"/><w:r><w:t xml:space="preserve">SRS-11: </w:t></w:r><w:r><w:t>The system shall display the error message “</w:t></w:r><w:r><w:t>The picture that must be image format</w:t></w:r><w:r><w:t>”</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>URS-10: The administrator view customer history</w:t></w:r><w:bookmarkEnd w:id="7"/></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>SRS-36</w:t></w:r><w:r><w:t>: The system shall provide interface for view a customer history.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-46</w:t></w:r><w:r><w:t>: The system shall provide search history page which receive customer username for searching.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-47</w:t></w:r><w:r><w:t>: The system query a list of history after searching by customer username</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>SRS-48</w:t></w:r><w:r><w:t>: The system shall provide search history page which receive date for searching.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Cordia New"/><w:szCs w:val="25"/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:tab/><w:t>SRS-49</w:t></w:r><w:r><w:t>: The system query a list of history after searching by</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New" w:hint="cs"/><w:szCs w:val="25"/><w:cs/><w:lang w:bidi="th-TH"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:szCs w:val="25"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>date.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:bookmarkStart w:id="8" w:name="_Toc379225282"/><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>URS-11: The customer save a product in cart</w:t></w:r><w:bookmarkEnd w:id="8"/></w:p><w:p><w:pPr><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/><w:t>SRS-50</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>: The system provide an interface that save product in</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t xml:space="preserve"> a</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t xml:space="preserve"> cart.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Cordia New"/><w:cs/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/><w:t>SRS-</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>51</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>: The system save a product in cart into a database.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:bookmarkStart w:id="9" w:name="_Toc379225283"/><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t>URS-12: The customer logout from system</w:t></w:r><w:bookmarkEnd w:id="9"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/><w:t xml:space="preserve">SRS- </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>52</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t xml:space="preserve"> The system provide an interface for logout.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:cs/><w:lang w:bidi="th-TH"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>SRS- 53</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>: The system not show a user status on top right after logout successful.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:bookmarkStart w:id="10" w:name="_Toc379225284"/><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:lastRenderedPageBreak/><w:t>URS-13: The administrator logout from system</w:t></w:r><w:bookmarkEnd w:id="10"/></w:p><w:p><w:pPr><w:rPr><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/><w:t>SRS- 52</w:t></w:r><w:bookmarkStart w:id="11" w:name="_GoBack"/><w:bookmarkEnd w:id="11"/><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:t xml:space="preserve"> the system provide an interface for logout.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:bidi="th-TH"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>SRS- 53</w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Cordia New"/><w:lang w:bidi="th-TH"/></w:rPr><w:t>: The system not show a user status on top right after logout successful.</w:t></w:r></w:p><w:p/>
'@

$d.Content.InsertXML($bodyXml)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
